$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.367.69"
$ws.Range("E2").Value = "  +2.07%  "

$ws.Range("D3").Value = "1.813.09"
$ws.Range("E3").Value = "  +2.56%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9976"
$ws.Range("E4").Value = "  -0.96%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.90"
$ws.Range("E5").Value = "  -0.23%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9955"
$ws.Range("E6").Value = "  -0.81%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4406"
$ws.Range("E7").Value = "  +16.36%  "

$ws.Range("E8").Value = "  +4.33%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.69"
$ws.Range("E9").Value = "  -0.89%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.159"
$ws.Range("E10").Value = "  +2.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07474"
$ws.Range("E11").Value = "  +2.33%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "23.00"
$ws.Range("E12").Value = "  -0.63%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.9964"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.295"
$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.299"
$ws.Range("E15").Value = "  -0.34%  "

$ws.Range("D16").Value = "1.816.89"
$ws.Range("E16").Value = "  +2.51%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001090"
$ws.Range("E17").Value = "  +2.59%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06668"
$ws.Range("E18").Value = "  +0.40%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "82.07"
$ws.Range("E19").Value = "  +0.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9959"
$ws.Range("E20").Value = "  -0.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.36"
$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.462"
$ws.Range("E22").Value = "  +1.53%  "

$ws.Range("D23").Value = "28.346.07"
$ws.Range("E23").Value = "  +1.92%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.08"
$ws.Range("E24").Value = "  +1.96%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.386"
$ws.Range("E25").Value = "  -0.31%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.465"
$ws.Range("E26").Value = "  +4.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.78"
$ws.Range("E27").Value = "  +2.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "155.31"
$ws.Range("E28").Value = "  +2.44%  "

$ws.Range("D29").Value = "2.020.66"
$ws.Range("E29").Value = "  +2.46%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.312"
$ws.Range("E30").Value = "  -11.69%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "132.89"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.061"
$ws.Range("E32").Value = "  +0.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.986"
$ws.Range("E33").Value = "  +1.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.09332"
$ws.Range("E34").Value = "  +5.98%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "12.34"
$ws.Range("E35").Value = "  -1.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6837"
$ws.Range("E36").Value = "  +1.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02388"
$ws.Range("E37").Value = "  +0.64%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06281"
$ws.Range("E38").Value = "  -0.37%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.225"
$ws.Range("E39").Value = "  +0.60%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2170"
$ws.Range("E40").Value = "  +1.72%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.495"
$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.223"
$ws.Range("E42").Value = "  +0.58%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.264"
$ws.Range("E43").Value = "  +1.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9948"
$ws.Range("E44").Value = "  -0.84%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.95"
$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6186"
$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.861"
$ws.Range("E47").Value = "  +0.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.11"
$ws.Range("E48").Value = "  -1.40%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.054"
$ws.Range("E49").Value = "  +0.91%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").Value = "  -0.99%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07117"
$ws.Range("E51").Value = "  -2.63%  "
